# Auto-generated edit script for cs-en-us-052pct.xlsx weekly update
# "New crime data collected" - updates volume/week header and crime stat table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/8/2024  Through  4/14/2024"

# --- Crime statistics table updates (rows 14-30) ---
$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E14").Value = -100
$ws.Range("G14").NumberFormat = '#,##0'
$ws.Range("G14").Value = 1
$ws.Range("H14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H14").Value = -100
$ws.Range("J14").NumberFormat = '#,##0'
$ws.Range("J14").Value = 1
$ws.Range("K14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K14").Value = 0
$ws.Range("C15").NumberFormat = 'General'
$ws.Range("C15").Formula = "'0"
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 33.333333333333
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = 54.545454545454
$ws.Range("M15").Value = 142.857142857143
$ws.Range("N15").Value = -5.555555555555
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 14
$ws.Range("E16").Value = -14.285714285714
$ws.Range("F16").Value = 31
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = 10.714285714285
$ws.Range("I16").Value = 147
$ws.Range("J16").Value = 110
$ws.Range("K16").Value = 33.636363636363
$ws.Range("L16").Value = 40
$ws.Range("M16").Value = 16.666666666666
$ws.Range("N16").Value = -66.438356164383
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -62.5
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = -5.128205128205
$ws.Range("I17").Value = 176
$ws.Range("J17").Value = 170
$ws.Range("K17").Value = 3.529411764705
$ws.Range("L17").Value = 13.548387096774
$ws.Range("M17").Value = 49.152542372881
$ws.Range("N17").Value = 0.571428571428
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 68
$ws.Range("K18").Value = -22.058823529411
$ws.Range("L18").Value = -15.873015873015
$ws.Range("M18").Value = -50.925925925925
$ws.Range("N18").Value = -91.666666666666
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = 25.490196078431
$ws.Range("I19").Value = 209
$ws.Range("J19").Value = 176
$ws.Range("K19").Value = 18.75
$ws.Range("L19").Value = 11.764705882352
$ws.Range("M19").Value = 39.333333333333
$ws.Range("N19").Value = -5.429864253393
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 26.315789473684
$ws.Range("I20").Value = 101
$ws.Range("J20").Value = 90
$ws.Range("K20").Value = 12.222222222222
$ws.Range("L20").Value = 29.487179487179
$ws.Range("M20").Value = 110.416666666667
$ws.Range("N20").Value = -76.993166287015
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 52
$ws.Range("E21").Value = -21.153846153846
$ws.Range("F21").Value = 170
$ws.Range("G21").Value = 161
$ws.Range("H21").Value = 5.590062111801
$ws.Range("I21").Value = 704
$ws.Range("J21").Value = 626
$ws.Range("K21").Value = 12.460063897763
$ws.Range("L21").Value = 17.333333333333
$ws.Range("M21").Value = 25.939177101967
$ws.Range("N21").Value = -63.730036063884
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 1
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 2
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 6
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 16
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = 23.076923076923
$ws.Range("L22").Value = 60
$ws.Range("M22").Value = 300
$ws.Range("C23").NumberFormat = '#,##0'
$ws.Range("C23").Value = 1
$ws.Range("F23").NumberFormat = '#,##0'
$ws.Range("F23").Value = 1
$ws.Range("I23").Value = 2
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -33.333333333333
$ws.Range("C24").Value = 48
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = 23.076923076923
$ws.Range("F24").Value = 154
$ws.Range("G24").Value = 241
$ws.Range("H24").Value = -36.09958506224
$ws.Range("I24").Value = 578
$ws.Range("J24").Value = 702
$ws.Range("K24").Value = -17.663817663817
$ws.Range("L24").Value = -37.378114842903
$ws.Range("M24").Value = 58.791208791208
$ws.Range("C25").Value = 30
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = 30.434782608695
$ws.Range("F25").Value = 87
$ws.Range("G25").Value = 181
$ws.Range("H25").Value = -51.933701657458
$ws.Range("I25").Value = 348
$ws.Range("J25").Value = 472
$ws.Range("K25").Value = -26.271186440678
$ws.Range("L25").Value = -53.036437246963
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 66
$ws.Range("H26").Value = -24.242424242424
$ws.Range("I26").Value = 235
$ws.Range("J26").Value = 208
$ws.Range("K26").Value = 12.980769230769
$ws.Range("L26").Value = 9.302325581395
$ws.Range("M26").Value = 1.293103448275
$ws.Range("C27").NumberFormat = 'General'
$ws.Range("C27").Formula = "'0"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = 4.347826086956
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 150
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 83.333333333333
$ws.Range("I28").Value = 33
$ws.Range("J28").Value = 25
$ws.Range("K28").Value = 32
$ws.Range("L28").Value = 32
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E29").Value = -100
$ws.Range("G29").NumberFormat = '#,##0'
$ws.Range("G29").Value = 1
$ws.Range("H29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 6
$ws.Range("K29").Value = 83.333333333333
$ws.Range("L29").Value = -15.384615384615
$ws.Range("D30").NumberFormat = '#,##0'
$ws.Range("D30").Value = 1
$ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E30").Value = -100
$ws.Range("G30").NumberFormat = '#,##0'
$ws.Range("G30").Value = 1
$ws.Range("H30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = 16.666666666666
$ws.Range("L30").Value = -36.363636363636

# --- Column E width adjustment (auto bestFit due to wider content) ---
$ws.Columns("E:E").AutoFit() | Out-Null
